$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes ---
# Row 8: Encoding Used -> TargetEncoder (was OneHotEncoder)
$ws.Range("B8").Value = "TargetEncoder"
$ws.Range("C8").Value = "TargetEncoder"
$ws.Range("D8").Value = "TargetEncoder"

# Row 9: Selection Method -> Top-25 Mutual Information (was Top-10)
$ws.Range("D9").Value = "Top-25 Mutual Information"

# Row 16: Notes bullet point Top-10 -> Top-25
$ws.Range("A16").Value = "• Mutual Information Top-25 selection used for stable dimensionality reduction."

# --- Selection change (cursor landed on C11 after edits) ---
$ws.Range("C11").Select()
